$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed TPM-derived metrics (ligand/receptor expression, specificity,
# and edge-weight columns G-J and M-T) for every data row 2-25 to match
# the updated NATMI TPM recompute.
$updates = @(
    @(2, 7, "0.053267"),
    @(2, 8, "0.159801"),
    @(2, 9, "0.0002370783698475801"),
    @(2, 10, "0.00023707836984758"),
    @(2, 13, "4.021066666666667"),
    @(2, 14, "12.0632"),
    @(2, 15, "0.0266528412086261"),
    @(2, 16, "0.0266528412086261"),
    @(2, 17, "0.2141901581333334"),
    @(2, 18, "1.9277114232"),
    @(2, 19, "6.318812145547482E-06"),
    @(2, 20, "6.318812145547481E-06"),
    @(3, 7, "0.053267"),
    @(3, 8, "0.159801"),
    @(3, 9, "0.0002370783698475801"),
    @(3, 10, "0.00023707836984758"),
    @(3, 15, "0.06472716901243536"),
    @(3, 16, "0.06472716901243537"),
    @(3, 17, "0.5201667791353334"),
    @(3, 18, "4.681501012218"),
    @(3, 19, "1.534541171431697E-05"),
    @(3, 20, "1.534541171431698E-05"),
    @(4, 7, "0.053267"),
    @(4, 8, "0.159801"),
    @(4, 9, "0.0002370783698475801"),
    @(4, 10, "0.00023707836984758"),
    @(4, 13, "16.16775866666667"),
    @(4, 14, "48.503276"),
    @(4, 15, "0.1071647749623786"),
    @(4, 16, "0.1071647749623786"),
    @(4, 17, "0.8612080008973334"),
    @(4, 18, "7.750872008076"),
    @(4, 19, "2.540645015316348E-05"),
    @(4, 20, "2.540645015316348E-05"),
    @(5, 7, "0.053267"),
    @(5, 8, "0.159801"),
    @(5, 9, "0.0002370783698475801"),
    @(5, 10, "0.00023707836984758"),
    @(5, 13, "5.277637333333334"),
    @(5, 14, "15.832912"),
    @(5, 15, "0.03498177012783927"),
    @(5, 16, "0.03498177012783928"),
    @(5, 17, "0.2811239078346667"),
    @(5, 18, "2.530115170512"),
    @(5, 19, "8.293421036290906E-06"),
    @(5, 20, "8.293421036290908E-06"),
    @(6, 7, "0.053267"),
    @(6, 8, "0.159801"),
    @(6, 9, "0.0002370783698475801"),
    @(6, 10, "0.00023707836984758"),
    @(6, 13, "108.7502723333333"),
    @(6, 14, "326.250817"),
    @(6, 15, "0.720829565926581"),
    @(6, 16, "0.7208295659265811"),
    @(6, 17, "5.792800756379666"),
    @(6, 18, "52.135206807417"),
    @(6, 19, "0.0001708930984278126"),
    @(6, 20, "0.0001708930984278126"),
    @(7, 7, "0.053267"),
    @(7, 8, "0.159801"),
    @(7, 9, "0.0002370783698475801"),
    @(7, 10, "0.00023707836984758"),
    @(7, 13, "6.886210666666667"),
    @(7, 14, "20.658632"),
    @(7, 15, "0.04564387876213955"),
    @(7, 16, "0.04564387876213956"),
    @(7, 17, "0.3668077835813334"),
    @(7, 18, "3.301270052232"),
    @(7, 19, "1.082117637044863E-05"),
    @(7, 20, "1.082117637044863E-05"),
    @(8, 7, "9.236317"),
    @(8, 8, "27.708951"),
    @(8, 9, "0.0411085846350553"),
    @(8, 10, "0.0411085846350553"),
    @(8, 13, "4.021066666666667"),
    @(8, 14, "12.0632"),
    @(8, 15, "0.0266528412086261"),
    @(8, 16, "0.0266528412086261"),
    @(8, 17, "37.13984641146667"),
    @(8, 18, "334.2586177032"),
    @(8, 19, "0.001095660578589496"),
    @(8, 20, "0.001095660578589496"),
    @(9, 7, "9.236317"),
    @(9, 8, "27.708951"),
    @(9, 9, "0.0411085846350553"),
    @(9, 10, "0.0411085846350553"),
    @(9, 15, "0.06472716901243536"),
    @(9, 16, "0.06472716901243537"),
    @(9, 17, "90.19515394076866"),
    @(9, 18, "811.756385466918"),
    @(9, 19, "0.002660842305535228"),
    @(9, 20, "0.002660842305535229"),
    @(10, 7, "9.236317"),
    @(10, 8, "27.708951"),
    @(10, 9, "0.0411085846350553"),
    @(10, 10, "0.0411085846350553"),
    @(10, 13, "16.16775866666667"),
    @(10, 14, "48.503276"),
    @(10, 15, "0.1071647749623786"),
    @(10, 16, "0.1071647749623786"),
    @(10, 17, "149.3305442248307"),
    @(10, 18, "1343.974898023476"),
    @(10, 19, "0.004405392221437596"),
    @(10, 20, "0.004405392221437596"),
    @(11, 7, "9.236317"),
    @(11, 8, "27.708951"),
    @(11, 9, "0.0411085846350553"),
    @(11, 10, "0.0411085846350553"),
    @(11, 13, "5.277637333333334"),
    @(11, 14, "15.832912"),
    @(11, 15, "0.03498177012783927"),
    @(11, 16, "0.03498177012783928"),
    @(11, 17, "48.74593142170134"),
    @(11, 18, "438.713382795312"),
    @(11, 19, "0.00143805105798433"),
    @(11, 20, "0.00143805105798433"),
    @(12, 7, "9.236317"),
    @(12, 8, "27.708951"),
    @(12, 9, "0.0411085846350553"),
    @(12, 10, "0.0411085846350553"),
    @(12, 13, "108.7502723333333"),
    @(12, 14, "326.250817"),
    @(12, 15, "0.720829565926581"),
    @(12, 16, "0.7208295659265811"),
    @(12, 17, "1004.451989106996"),
    @(12, 18, "9040.067901962966"),
    @(12, 19, "0.02963228321834303"),
    @(12, 20, "0.02963228321834304"),
    @(13, 7, "9.236317"),
    @(13, 8, "27.708951"),
    @(13, 9, "0.0411085846350553"),
    @(13, 10, "0.0411085846350553"),
    @(13, 13, "6.886210666666667"),
    @(13, 14, "20.658632"),
    @(13, 15, "0.04564387876213955"),
    @(13, 16, "0.04564387876213956"),
    @(13, 17, "63.60322464611467"),
    @(13, 18, "572.429021815032"),
    @(13, 19, "0.001876355253165617"),
    @(13, 20, "0.001876355253165617"),
    @(14, 7, "212.661977"),
    @(14, 8, "637.9859310000001"),
    @(14, 9, "0.9465063704680865"),
    @(14, 10, "0.9465063704680865"),
    @(14, 13, "4.021066666666667"),
    @(14, 14, "12.0632"),
    @(14, 15, "0.0266528412086261"),
    @(14, 16, "0.0266528412086261"),
    @(14, 17, "855.1279869821334"),
    @(14, 18, "7696.151882839201"),
    @(14, 19, "0.02522708399503894"),
    @(14, 20, "0.02522708399503894"),
    @(15, 7, "212.661977"),
    @(15, 8, "637.9859310000001"),
    @(15, 9, "0.9465063704680865"),
    @(15, 10, "0.9465063704680865"),
    @(15, 15, "0.06472716901243536"),
    @(15, 16, "0.06472716901243537"),
    @(15, 17, "2076.702191237395"),
    @(15, 18, "18690.31972113656"),
    @(15, 19, "0.06126467781263459"),
    @(15, 20, "0.0612646778126346"),
    @(16, 7, "212.661977"),
    @(16, 8, "637.9859310000001"),
    @(16, 9, "0.9465063704680865"),
    @(16, 10, "0.9465063704680865"),
    @(16, 13, "16.16775866666667"),
    @(16, 14, "48.503276"),
    @(16, 15, "0.1071647749623786"),
    @(16, 16, "0.1071647749623786"),
    @(16, 17, "3438.267521712218"),
    @(16, 18, "30944.40769540996"),
    @(16, 19, "0.1014321421916702"),
    @(16, 20, "0.1014321421916702"),
    @(17, 7, "212.661977"),
    @(17, 8, "637.9859310000001"),
    @(17, 9, "0.9465063704680865"),
    @(17, 10, "0.9465063704680865"),
    @(17, 13, "5.277637333333334"),
    @(17, 14, "15.832912"),
    @(17, 15, "0.03498177012783927"),
    @(17, 16, "0.03498177012783928"),
    @(17, 17, "1122.352789195675"),
    @(17, 18, "10101.17510276107"),
    @(17, 19, "0.03311046827625008"),
    @(17, 20, "0.03311046827625008"),
    @(18, 7, "212.661977"),
    @(18, 8, "637.9859310000001"),
    @(18, 9, "0.9465063704680865"),
    @(18, 10, "0.9465063704680865"),
    @(18, 13, "108.7502723333333"),
    @(18, 14, "326.250817"),
    @(18, 15, "0.720829565926581"),
    @(18, 16, "0.7208295659265811"),
    @(18, 17, "23127.04791369507"),
    @(18, 18, "208143.4312232556"),
    @(18, 19, "0.6822697761712544"),
    @(18, 20, "0.6822697761712545"),
    @(19, 7, "212.661977"),
    @(19, 8, "637.9859310000001"),
    @(19, 9, "0.9465063704680865"),
    @(19, 10, "0.9465063704680865"),
    @(19, 13, "6.886210666666667"),
    @(19, 14, "20.658632"),
    @(19, 15, "0.04564387876213955"),
    @(19, 16, "0.04564387876213956"),
    @(19, 17, "1464.435174411821"),
    @(19, 18, "13179.91656970639"),
    @(19, 19, "0.04320222202123809"),
    @(19, 20, "0.04320222202123809"),
    @(20, 7, "2.729417"),
    @(20, 8, "8.188250999999999"),
    @(20, 9, "0.01214796652701058"),
    @(20, 10, "0.01214796652701058"),
    @(20, 13, "4.021066666666667"),
    @(20, 14, "12.0632"),
    @(20, 15, "0.0266528412086261"),
    @(20, 16, "0.0266528412086261"),
    @(20, 17, "10.97516771813333"),
    @(20, 18, "98.77650946319999"),
    @(20, 19, "0.000323777822852118"),
    @(20, 20, "0.000323777822852118"),
    @(21, 7, "2.729417"),
    @(21, 8, "8.188250999999999"),
    @(21, 9, "0.01214796652701058"),
    @(21, 10, "0.01214796652701058"),
    @(21, 15, "0.06472716901243536"),
    @(21, 16, "0.06472716901243537"),
    @(21, 17, "26.65350122603533"),
    @(21, 18, "239.881511034318"),
    @(21, 19, "0.0007863034825512209"),
    @(21, 20, "0.0007863034825512211"),
    @(22, 7, "2.729417"),
    @(22, 8, "8.188250999999999"),
    @(22, 9, "0.01214796652701058"),
    @(22, 10, "0.01214796652701058"),
    @(22, 13, "16.16775866666667"),
    @(22, 14, "48.503276"),
    @(22, 15, "0.1071647749623786"),
    @(22, 16, "0.1071647749623786"),
    @(22, 17, "44.12855535669733"),
    @(22, 18, "397.156998210276"),
    @(22, 19, "0.001301834099117596"),
    @(22, 20, "0.001301834099117596"),
    @(23, 7, "2.729417"),
    @(23, 8, "8.188250999999999"),
    @(23, 9, "0.01214796652701058"),
    @(23, 10, "0.01214796652701058"),
    @(23, 13, "5.277637333333334"),
    @(23, 14, "15.832912"),
    @(23, 15, "0.03498177012783927"),
    @(23, 16, "0.03498177012783928"),
    @(23, 17, "14.40487305743467"),
    @(23, 18, "129.643857516912"),
    @(23, 19, "0.0004249573725685699"),
    @(23, 20, "0.00042495737256857"),
    @(24, 7, "2.729417"),
    @(24, 8, "8.188250999999999"),
    @(24, 9, "0.01214796652701058"),
    @(24, 10, "0.01214796652701058"),
    @(24, 13, "108.7502723333333"),
    @(24, 14, "326.250817"),
    @(24, 15, "0.720829565926581"),
    @(24, 16, "0.7208295659265811"),
    @(24, 17, "296.8248420612296"),
    @(24, 18, "2671.423578551067"),
    @(24, 19, "0.008756613438555669"),
    @(24, 20, "0.008756613438555669"),
    @(25, 7, "2.729417"),
    @(25, 8, "8.188250999999999"),
    @(25, 9, "0.01214796652701058"),
    @(25, 10, "0.01214796652701058"),
    @(25, 13, "6.886210666666667"),
    @(25, 14, "20.658632"),
    @(25, 15, "0.04564387876213955"),
    @(25, 16, "0.04564387876213956"),
    @(25, 17, "18.79534045918133"),
    @(25, 18, "169.158064132632"),
    @(25, 19, "0.0005544803113654001"),
    @(25, 20, "0.0005544803113654003")

)

foreach ($u in $updates) {
    $row = [int]$u[0]
    $col = [int]$u[1]
    $val = [double]$u[2]
    $ws.Cells.Item($row, $col).Value2 = $val
}
